$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the forward-slash path separators with backslashes, and swap the
# "CS+4 / UCS4 / laugh1" condition block (rows 17-31 originally split into a
# CS+4/black/silent sub-block and a CS+4/UCS4/laugh1 sub-block) for the
# "huang" folder's CS-3 condition (rows 17-31 now uniformly reference
# CS-3 / black / silent).

# Column A (CSName)
$ws.Range("A2:A16").Value = "Condition\CS+3.BMP"
$ws.Range("A17:A31").Value = "Condition\CS-3.BMP"

# Column B (UCSName)
$ws.Range("B2:B4").Value = "Condition\black.PNG"
$ws.Range("B5:B16").Value = "Condition\UCS3.BMP"
$ws.Range("B17:B31").Value = "Condition\black.PNG"

# Column C (SoundName)
$ws.Range("C2:C4").Value = "Sound\silent.wav"
$ws.Range("C5:C16").Value = "Sound\scream1.wav"
$ws.Range("C17:C31").Value = "Sound\silent.wav"

# Column D (TrgCol) values are unchanged (-0.1 / 0.9 / 0.4), left as-is.
